$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before writing, so numeric-looking
# strings like "1.024" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.548.84'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '1.877.47'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").Value = '1.024'
$ws.Range("E4").Value = '  +1.78%  '
$ws.Range("D5").Value = '318.37'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").Value = '1.023'
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("D7").Value = '0.5149'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '0.3956'
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("D9").Value = '0.08341'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = '1.117'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '42.12'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '6.263'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '20.53'
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").Value = '1.859.63'
$ws.Range("E14").Value = '  -1.27%  '
$ws.Range("D15").Value = '1.024'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '7.244'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '0.00001111'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '91.43'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = '0.06792'
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").Value = '17.73'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '1.023'
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").Value = '5.992'
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").Value = '28.585.02'
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("D24").Value = '11.16'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").Value = '2.278'
$ws.Range("E25").Value = '  +0.51%  '
$ws.Range("D26").Value = '2.071.19'
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("D27").Value = '162.06'
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").Value = '20.81'
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D29").Value = '2.375'
$ws.Range("E29").Value = '  -4.32%  '
$ws.Range("D30").Value = '127.71'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("D32").Value = '1.036'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").Value = '5.837'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = '3.649'
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").Value = '0.02439'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("D36").Value = '0.06520'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("D37").Value = '9.200'
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("D38").Value = '0.2183'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D40").Value = '1.187'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").Value = '0.6444'
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").Value = '5.007'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '11.20'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("D44").Value = '0.6051'
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").Value = '13.02'
$ws.Range("E45").Value = '  -0.61%  '
$ws.Range("D46").Value = '3.715'
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("D47").Value = '1.238'
$ws.Range("E47").Value = '  -3.62%  '
$ws.Range("D48").Value = '1.996'
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("D49").Value = '1.213'
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("D50").Value = '122.07'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").Value = '0.06857'
$ws.Range("E51").Value = '  -0.69%  '

# Restore the Normal style on column D so no stray cell-style / number-format
# metadata is left behind (matches original inlineStr plain-text cells).
$ws.Range("D2:D51").Style = "Normal"
